# Save User 버튼을 individual로 수정
$wb = $excel.ActiveWorkbook

# --- Sheet "Group": add row 5 (A5 = "TEST", B5 = empty inline string) ---
# Note: B5 is an empty (but present) cell in the source; Excel's COM model
# treats Value = "" as "clear the cell", which is the closest reachable
# approximation (openpyxl / most readers surface both as None).
$wsGroup = $wb.Worksheets.Item("Group")
$wsGroup.Range("A5").Value = "TEST"
$wsGroup.Range("B5").Value = ""

# --- Sheet "Users": fix row 13 types and add rows 14-16 ---
$wsUsers = $wb.Worksheets.Item("Users")

# Row 13: convert A13, D13, E13, F13 to numeric values
$wsUsers.Cells.Item(13, 1).Value = 12
$wsUsers.Cells.Item(13, 2).Value = "park"
$wsUsers.Cells.Item(13, 3).Value = "male"
$wsUsers.Cells.Item(13, 4).Value = 53
$wsUsers.Cells.Item(13, 5).Value = 170
$wsUsers.Cells.Item(13, 6).Value = 72
$wsUsers.Cells.Item(13, 7).Value = "test1"

# Row 14
$wsUsers.Cells.Item(14, 1).Value = "testuser01"
$wsUsers.Cells.Item(14, 2).Value = "Park"
$wsUsers.Cells.Item(14, 3).Value = "male"
$wsUsers.Cells.Item(14, 4).Value = 53
$wsUsers.Cells.Item(14, 5).Value = 170
$wsUsers.Cells.Item(14, 6).Value = 72
$wsUsers.Cells.Item(14, 7).Value = "TEST"

# Row 15
$wsUsers.Cells.Item(15, 1).Value = "testuser02"
$wsUsers.Cells.Item(15, 2).Value = "Lee"
$wsUsers.Cells.Item(15, 3).Value = "female"
$wsUsers.Cells.Item(15, 4).Value = 48
$wsUsers.Cells.Item(15, 5).Value = 148
$wsUsers.Cells.Item(15, 6).Value = 56
$wsUsers.Cells.Item(15, 7).Value = "TEST"

# Row 16 (age/height/weight stored as text, not numbers)
$wsUsers.Cells.Item(16, 1).Value = "testuser03"
$wsUsers.Cells.Item(16, 2).Value = "Kim"
$wsUsers.Cells.Item(16, 3).Value = "female"
$wsUsers.Cells.Item(16, 4).Value = "'28"
$wsUsers.Cells.Item(16, 5).Value = "'163.8"
$wsUsers.Cells.Item(16, 6).Value = "'53"
$wsUsers.Cells.Item(16, 7).Value = "TEST"
